$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update Price (D) and Volume(1h) (E) columns,
# and fix two pairs of mis-ordered rows (14/15 and 44/45).

$ws.Range('D2').Value = '44.569.22'
$ws.Range('E2').Value = '  +1.43%  '
$ws.Range('D3').Value = '2.248.57'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.64'
$ws.Range('E5').Value = '  +1.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.29'
$ws.Range('E6').Value = '  +0.70%  '
$ws.Range('E7').Value = '  +1.18%  '
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('E9').Value = '  +1.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.14'
$ws.Range('E10').Value = '  +2.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0811'
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.29'
$ws.Range('E12').Value = '  +2.09%  '
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '2.320.80'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.843'
$ws.Range('E15').Value = '  +4.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.64'
$ws.Range('E16').Value = '  +2.10%  '
$ws.Range('D17').Value = '44.236.78'
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D18').Value = '0.0₃0966'
$ws.Range('E18').Value = '  +1.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.41'
$ws.Range('E19').Value = '  +4.42%  '
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '65.88'
$ws.Range('E21').Value = '  +2.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '239.82'
$ws.Range('E22').Value = '  +1.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.00'
$ws.Range('E23').Value = '  +3.47%  '
$ws.Range('E24').Value = '  +4.56%  '
$ws.Range('E25').Value = '  -0.36%  '
$ws.Range('E26').Value = '  +5.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.85'
$ws.Range('E27').Value = '  +0.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.63'
$ws.Range('E28').Value = '  +4.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.04'
$ws.Range('E29').Value = '  +3.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.14'
$ws.Range('E30').Value = '  +0.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '152.83'
$ws.Range('E31').Value = '  +0.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0805'
$ws.Range('E32').Value = '  +0.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.65'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.16'
$ws.Range('E34').Value = '  -2.45%  '
$ws.Range('E35').Value = '  +1.42%  '
$ws.Range('E36').Value = '  +2.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.82'
$ws.Range('E37').Value = '  +3.53%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.45'
$ws.Range('E38').Value = '  +4.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.47'
$ws.Range('E39').Value = '  -1.46%  '
$ws.Range('E40').Value = '  +1.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0303'
$ws.Range('E41').Value = '  +2.32%  '
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('D43').Value = '1.756.30'
$ws.Range('E43').Value = '  +1.86%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.194'
$ws.Range('E44').Value = '  +5.41%  '
$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '81.14'
$ws.Range('E45').Value = '  -3.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '71.27'
$ws.Range('E46').Value = '  +4.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '99.89'
$ws.Range('E47').Value = '  +0.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.68'
$ws.Range('E48').Value = '  +3.79%  '
$ws.Range('E49').Value = '  +2.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.88'
$ws.Range('E50').Value = '  +0.71%  '
$ws.Range('E51').Value = '  +4.79%  '
